# Week 17 update:
#  - Rushing sheet: D.Henderson released from the roster (row removed);
#    S.Michel's season rushing totals updated with Week 17 stats.
#  - Receiving sheet: Week 17 receiving stats added for several players.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet -------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# D.Henderson (row 4) is no longer on the roster - remove his row entirely.
# This shifts everyone below him up by one row.
$rushing.Rows.Item(4).Delete()

# S.Michel (now row 4 after the shift) gets his Week 17 totals added.
$rushing.Range("C4").Value = 86
$rushing.Range("D4").Value = 46
$rushing.Range("E4").Value = 31
$rushing.Range("F4").Value = 29

# Re-sequence the roster-index column (A) now that a player was removed.
$rushing.Range("A2").Value = 0
$rushing.Range("A3").Value = 1
$rushing.Range("A4").Value = 2
$rushing.Range("A5").Value = 3
$rushing.Range("A6").Value = 4
$rushing.Range("A7").Value = 5
$rushing.Range("A8").Value = 6
$rushing.Range("A9").Value = 7

# --- Receiving sheet -------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# S.Michel (row 2)
$receiving.Range("C2").Value = 43
$receiving.Range("D2").Value = 30

# C.Kupp (row 3)
$receiving.Range("C3").Value = 128
$receiving.Range("D3").Value = 99
$receiving.Range("E3").Value = 44
$receiving.Range("F3").Value = 26

# V.Jefferson (row 4)
$receiving.Range("C4").Value = 59
$receiving.Range("D4").Value = 36
$receiving.Range("E4").Value = 26
$receiving.Range("F4").Value = 10

# B.Skowronek (row 5)
$receiving.Range("C5").Value = 24
$receiving.Range("D5").Value = 9

# O.Beckham (row 6)
$receiving.Range("C6").Value = 64
$receiving.Range("D6").Value = 54
$receiving.Range("E6").Value = 27

# T.Higbee (row 8)
$receiving.Range("C8").Value = 69
$receiving.Range("D8").Value = 52

# Re-sequence the roster-index column (A) - it used to skip index 1
# (D.Henderson had no receiving stats); with him off the roster entirely
# the remaining players' indices close the gap.
$receiving.Range("A2").Value = 0
$receiving.Range("A3").Value = 1
$receiving.Range("A4").Value = 2
$receiving.Range("A5").Value = 3
$receiving.Range("A6").Value = 4
$receiving.Range("A7").Value = 5
$receiving.Range("A8").Value = 6
$receiving.Range("A9").Value = 7
$receiving.Range("A10").Value = 8
